# The mistake of Worldview-2 (row 12) and Worldview-3 (row 13) was fixed:
# some Yellow (I/J) and Cirrus/SWIR (Y:AN) band values had been entered
# on the Worldview-2 row instead of the Worldview-3 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Folha3" - the active tab in this workbook

# --- Yellow_min / Yellow_max (columns I and J) ---
# Row 12 (Worldview-2) keeps only Yellow_min, equal to Green_max (H12)
$ws.Range("I12").Value = 0.57999999999999996
$ws.Range("J12").ClearContents()

# Row 13 (Worldview-3) gets the Yellow_min/_max values that used to sit on row 12
$ws.Range("I13").Value = 0.58499999999999996
$ws.Range("J13").Value = 0.625

# --- Cirrus_min/_max through SWIR7_min/_max (columns Y through AN) ---
# These belong to Worldview-3 (row 13), not Worldview-2 (row 12)
$cols   = @("Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN")
$values = @(1.1950000000000001,1.2250000000000001,1.55,1.59,1.64,1.68,1.71,1.75,2.145,2.1850000000000001,2.1850000000000001,2.2250000000000001,2.2349999999999999,2.2850000000000001,2.2949999999999999,2.3650000000000002)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").ClearContents()
    $ws.Range($cols[$i] + "13").Value = $values[$i]
}

# --- Selection left on row 12 (Worldview-2) after the edit ---
$ws.Range("C12:XFD12").Select()
